$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new export timestamp
$ws.Name = "IClientBalance-20240617-091042-"

# Update the "Dt. Referencia" column (G) for every data row (2-257)
# from 45457 (2024-06-14) to 45460 (2024-06-17)
$ws.Range("G2:G257").Value = 45460

# Row 103 also had its projected value (E) and total value (H) updated
$ws.Range("E103").Value = -12057.07
$ws.Range("H103").Value = 3015.66
